$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column "Seguridad" before the existing "Empleado por
#    contrato" column (old column G / 7).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).Insert()

$ws.Cells.Item(1, 7).Value = "Seguridad"

For ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 7).Value = "No"
}

# ---------------------------------------------------------------------------
# 2) Insert 12 new "seguridad" related columns right before the
#    "Seguro Social (9.75%)" column (now column S / 19 after step 1).
# ---------------------------------------------------------------------------
$ws.Range("S1:AD1").EntireColumn.Insert()

$securityHeaders = @(
    "Horas Turno Seguridad",
    "Hora Cambio Turno Seguridad",
    "Margen Salida Seguridad (min)",
    "Tolerancia Turno Seguridad (min)",
    "Horas Reales Seguridad (prom)",
    "Dif Turno Seguridad (min, prom)",
    "Alerta Seguridad",
    "Turnos Seguridad Día",
    "Turnos Seguridad Noche",
    "Total Turnos Seguridad",
    "Empleados Seguridad Turno Día",
    "Empleados Seguridad Turno Noche"
)

$col = 19
Foreach ($h in $securityHeaders) {
    $ws.Cells.Item(1, $col).Value = $h
    For ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, $col).Value = ""
    }
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 3) Update the loan balance for row 5 (now column AJ / 36):
#    "Total Saldo Préstamo" goes from 40 to 39.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 36).Value = 39
